# Updated 2018 lures test
# Update the ANOVA result for the "pi_md" row of the "script3" experiment
# (cell F7) to reflect the corrected statistic.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stats-results")

$ws.Range("F7").Value = "F = 10.86, df = 4, 25; P < 0.001"

# Mirror the author's on-screen state: scrolled so column B is the
# left-most visible column, with F7 as the active selection.
$ws.Activate()
$ws.Range("F7").Select()
$excel.ActiveWindow.ScrollColumn = 2
